$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new probability value (stored as text, matching original string formatting)
$updates = @{
    2  = "1.0000100135803223"
    4  = "0.9999911785125732"
    5  = "0.9993734955787659"
    6  = "0.9995195865631104"
    7  = "0.9763948917388916"
    8  = "0.991457462310791"
    9  = "0.9991353154182434"
    11 = "0.9991846680641174"
    12 = "0.952233612537384"
    13 = "0.9980811476707458"
    14 = "0.9981659054756165"
    15 = "0.9970064163208008"
    16 = "0.9297232031822205"
    17 = "0.9919484257698059"
    18 = "0.9998220801353455"
    19 = "0.998589277267456"
    20 = "0.9988031387329102"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("C$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
